$d = $word.ActiveDocument

# 1. After the "See code." paragraph, insert the new paragraphs:
#      NLP tasks:
#      Text generation            (bulleted list, same list as "In this course:")
#      Sentence similarity        (bulleted list)
#      Summarization               (bulleted list)
#      Machines translation       (bulleted list)
#    followed by one blank trailing paragraph (a second one is added afterwards).
$find = "See code."
$replace = "See code.^pNLP tasks:^pText generation^pSentence similarity^pSummarization^pMachines translation^p"
$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null

# 2. Add the second trailing blank paragraph at the very end of the document.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter() | Out-Null

# 3. Turn the four task paragraphs into bulleted list items that reuse the
#    same list (numId) that is already used by the "In this course:" bullets.
$listTemplate = $d.Paragraphs.Item(9).Range.ListFormat.ListTemplate
$itemTexts = @("Text generation", "Sentence similarity", "Summarization", "Machines translation")

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($itemTexts -contains $text) {
        $para.Style = "ListParagraph"
        $para.Range.ListFormat.ApplyListTemplate($listTemplate, $true) | Out-Null
    }
}
